$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = 0.00153
    $ws.Range("E$row").ClearContents()

    $ws.Range("G$row").Value = 0.0141267787839586
    $ws.Range("H$row").Value = 0.0141267787839586
    $ws.Range("I$row").Value = -0.01604139715394567
    $ws.Range("J$row").Value = -0.01604139715394567
    $ws.Range("K$row").Value = -3.43
    $ws.Range("L$row").Value = -0.0088745148771022
    $ws.Range("M$row").Value = 2.13
    $ws.Range("N$row").Value = 0.03317757009345794
    $ws.Range("O$row").Value = -0.6209912536443148
    $ws.Range("P$row").Value = 2.13
    $ws.Range("Q$row").Value = 0.03317757009345794
    $ws.Range("R$row").Value = -0.6209912536443148

    $ws.Range("U$row").Value = 2.13
    $ws.Range("V$row").Value = 0.03317757009345794
    $ws.Range("W$row").Value = -0.03235849056603773
    $ws.Range("X$row").Value = 0.2323193234261386
    $ws.Range("Y$row").Value = -0.2646778139921763
    $ws.Range("Z$row").Value = 2.342708207055401
    $ws.Range("AA$row").Value = -0.03758031276518366
    $ws.Range("AB$row").Value = 0.1357841962197665
    $ws.Range("AC$row").Value = -0.1733645089849501
    $ws.Range("AD$row").Value = 93.2
    $ws.Range("AE$row").Value = 0
    $ws.Range("AF$row").Value = 93.2
    $ws.Range("AG$row").Value = 91.07000000000001
    $ws.Range("AH$row").Value = 0.5921219822109276
    $ws.Range("AI$row").Value = 0.4801648634724369
    $ws.Range("AJ$row").Value = 0.5865266954337606
    $ws.Range("AK$row").Value = 0.4743970412043548
    $ws.Range("AL$row").Value = 1.88
    $ws.Range("AM$row").Value = -4.81
    $ws.Range("AN$row").Value = -23.01234567901235
    $ws.Range("AO$row").Value = -3.297872340425532
    $ws.Range("AP$row").Value = -22.48641975308642
    $ws.Range("AQ$row").Value = 1.288981288981289
}
